$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-18.94012429983433,2.076026314430466,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,1.358480904272813,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,2.069549209202879,-18.94012429983433,-18.94012429983433,3.533791726427909,-18.94012429983433,1.488245074536467,-18.94012429983433,0.9339039748352248,-18.94012429983433),
    @(-18.94012429983433,1.428651763647877,-18.94012429983433,-18.94012429983433,-18.94012429983433,2.827716245891096,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(2.475015525492851,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,1.937846509522344,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(3.852021966307857,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,1.688507239880275,-18.94012429983433,2.214818007931721),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,2.906789626885464,-18.94012429983433,2.858661123360293,-18.94012429983433,-18.94012429983433,-18.94012429983433,1.922053698189585),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,2.521455564260975,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,1.648152052025548,1.817588796463223),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,1.968629659017553),
    @(-18.94012429983433,-18.94012429983433,4.321925369971993,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,1.947843707162543,-18.94012429983433),
    @(-18.94012429983433,2.016360337950242,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,2.076543228790617,2.062685909637461,2.521355647055241,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,2.026238214945414,2.032870205818259,2.421162524172947,-18.94012429983433),
    @(-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,-18.94012429983433,1.600916526496476,1.824440381114828,-18.94012429983433,-18.94012429983433),
    @(-18.94012429983433,0.9289408334568783,-18.94012429983433,-18.94012429983433,3.073499953209746,-18.94012429983433,1.651899240388024,1.263683697253407,-18.94012429983433,2.045876718116762),
    @(-18.94012429983433,1.558036366822252,-18.94012429983433,1.549693288005597,-18.94012429983433,2.497498255183933,1.447473414156256,-18.94012429983433,-18.94012429983433,-18.94012429983433),
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
